$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to keep numeric-looking values as text
# (matches the workbook's original inlineStr/text cell formatting).

$ws.Range("D2").Value = "62.341.39"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "3.441.62"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'410.34"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'129.87"
$ws.Range("E6").Value = "  -4.04%  "
$ws.Range("D7").Value = "'0.633"
$ws.Range("E7").Value = "  +7.29%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.759"
$ws.Range("E9").Value = "  +12.72%  "
$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = "  +18.73%  "
$ws.Range("D11").Value = "'43.29"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "'8.94"
$ws.Range("E13").Value = "  +7.60%  "
$ws.Range("D14").Value = "'20.51"
$ws.Range("E14").Value = "  +4.99%  "
$ws.Range("D15").Value = "'0.0000196"
$ws.Range("E15").Value = "  +54.57%  "
$ws.Range("D16").Value = "3.459.73"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("E17").Value = "  +3.84%  "
$ws.Range("D18").Value = "62.311.61"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").Value = "'11.57"
$ws.Range("E19").Value = "  +5.73%  "
$ws.Range("D20").Value = "'396.13"
$ws.Range("E20").Value = "  +26.86%  "
$ws.Range("D21").Value = "'89.30"
$ws.Range("E21").Value = "  +6.85%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'13.34"
$ws.Range("E23").Value = "  +5.18%  "
$ws.Range("D24").Value = "'3.23"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").Value = "'32.00"
$ws.Range("E25").Value = "  +8.87%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").Value = "'8.52"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("E28").Value = "  +3.53%  "
$ws.Range("E29").Value = "  +11.22%  "
$ws.Range("E30").Value = "  +7.34%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "'0.172"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "'11.84"
$ws.Range("E33").Value = "  +5.06%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "'0.0494"
$ws.Range("E35").Value = "  +3.14%  "
$ws.Range("D36").Value = "'52.38"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "'0.132"
$ws.Range("E40").Value = "  +7.70%  "
$ws.Range("E41").Value = "  +7.83%  "
$ws.Range("D42").Value = "'140.83"
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'16.79"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("D47").Value = "'22.28"
$ws.Range("E47").Value = "  +4.88%  "
$ws.Range("D48").Value = "2.125.79"
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "'1.94"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'2.28"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").Value = "'0.0370"
$ws.Range("E51").Value = "  +8.42%  "
